$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric value changes ---
$ws.Cells.Item(6, 12).Value = 165
$ws.Cells.Item(7, 12).Value = 3
$ws.Cells.Item(21, 15).Value = 14
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(22, 15).Value = 14
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(23, 15).Value = 14
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(24, 15).Value = 13
$ws.Cells.Item(24, 16).Value = 1
$ws.Cells.Item(25, 15).Value = 14
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(26, 15).Value = 14
$ws.Cells.Item(26, 16).Value = 0

# --- Percent-like text changes (use Formula-as-text-literal to avoid Excel auto percent-number conversion and preserve style) ---
$ws.Cells.Item(9, 12).Formula = '="51.9%"'
$ws.Cells.Item(10, 12).Formula = '="75.1%"'
$ws.Cells.Item(21, 18).Formula = '="51.9%"'
$ws.Cells.Item(21, 19).Formula = '="77.6%"'
$ws.Cells.Item(22, 18).Formula = '="51.9%"'
$ws.Cells.Item(22, 19).Formula = '="76.0%"'
$ws.Cells.Item(23, 18).Formula = '="51.9%"'
$ws.Cells.Item(23, 19).Formula = '="80.6%"'
$ws.Cells.Item(24, 18).Formula = '="48.1%"'
$ws.Cells.Item(24, 19).Formula = '="70.9%"'
$ws.Cells.Item(25, 18).Formula = '="51.9%"'
$ws.Cells.Item(25, 19).Formula = '="69.2%"'
$ws.Cells.Item(26, 18).Formula = '="51.9%"'
$ws.Cells.Item(26, 19).Formula = '="62.8%"'

# --- G column swap: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$ws.Cells.Item(8, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(9, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(10, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(34, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(35, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(36, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(60, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(61, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(62, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(86, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(87, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(88, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(112, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(113, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(114, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(138, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(139, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(140, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(164, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(167, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(170, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(191, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(194, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(197, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(218, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(221, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(224, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(245, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(248, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(251, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(272, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(275, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(278, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(299, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(302, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(305, 7).Value = "System, dnasr281@gmail.com"

# --- Restructure rows: Not Recorded -> Recorded (style 4 -> 2, fill formats copied; content updated) ---
$src = $ws.Range("A2:I2")
$dst = $ws.Range("A171:I171")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(171, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(171, 8).Value = "20/23"
$ws.Cells.Item(171, 9).Value = "Recorded"

$src = $ws.Range("A2:I2")
$dst = $ws.Range("A198:I198")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(198, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(198, 8).Value = "22/30"
$ws.Cells.Item(198, 9).Value = "Recorded"

$src = $ws.Range("A2:I2")
$dst = $ws.Range("A225:I225")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(225, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(225, 8).Value = "19/25"
$ws.Cells.Item(225, 9).Value = "Recorded"

$src = $ws.Range("A2:I2")
$dst = $ws.Range("A252:I252")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(252, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(252, 8).Value = "20/28"
$ws.Cells.Item(252, 9).Value = "Recorded"

$src = $ws.Range("A2:I2")
$dst = $ws.Range("A279:I279")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(279, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(279, 8).Value = "15/26"
$ws.Cells.Item(279, 9).Value = "Recorded"

$src = $ws.Range("A2:I2")
$dst = $ws.Range("A306:I306")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(306, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(306, 8).Value = "22/29"
$ws.Cells.Item(306, 9).Value = "Recorded"
